# Remove "keywords" from Performance model
# Column O on Sheet1 (row 9) holds a duplicate "-" value (a leftover
# "keywords" field). Delete the whole column so everything to its right
# (similes SMILES, InChIKey, the two numeric fields, solvent and the
# "carbazole" donor group) shifts one column to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("O:O").Delete()

# Leave the selection on column N, mirroring where the cursor ended up
# after the deletion in the authoring session.
$ws.Columns("N:N").Select()
